$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 value changes to 0, C2 becomes a plain value (no formula), D2 becomes empty
$ws.Range("B2").Value = 0
$ws.Range("C2").Formula = 0
$ws.Range("D2").Value = ""

# Row 2 reverts to default (non-custom) row height
$ws.Rows(2).AutoFit()

# Row 3: B3 value changes to 5805, C3 formula recalculates automatically
$ws.Range("B3").Value = 5805

# Update selection to C3
$ws.Range("C3").Select()

